# Upgrade left table: add year 2023 (column K) with new data rows, matching
# the formatting of the existing year columns and closing the table with a
# right-hand border.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Copy the formatting of the last existing data column (J) onto the new
#    column K for each of the table rows, then fill in the 2023 values.
# ---------------------------------------------------------------------

# Row 3 (year headers)
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("K3").Value = 2023

# Row 4 (Employed persons)
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)
$ws.Range("K4").Value = 828

# Row 5 (Of which: Women)
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)
$ws.Range("K5").Value = 228

# Row 6 (Men)
$ws.Range("J6").Copy()
$ws.Range("K6").PasteSpecial(-4122)
$ws.Range("K6").Value = 600

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Column K now sits at the right-hand edge of the table, so it needs a
#    closing thin right border, matching whatever top/bottom border the
#    rest of the row already carries.
# ---------------------------------------------------------------------

$ws.Range("K3:K6").Borders.Item(10).LineStyle = 1   # xlEdgeRight / xlContinuous
$ws.Range("K3:K6").Borders.Item(10).Weight = 2      # xlThin

# ---------------------------------------------------------------------
# 3. Widen the custom-width column block so it now spans through column M
#    (same width as the existing year columns), matching the workbook's
#    updated column formatting.
# ---------------------------------------------------------------------

$ws.Range("K1:M1").EntireColumn.ColumnWidth = 7.8

# ---------------------------------------------------------------------
# 4. Row spans / sheet dimension are maintained automatically by the
#    engine when the new cells are written above.
# ---------------------------------------------------------------------
